# Fruta / hortaliza, semanal
# Insert two new weekly-update rows into the daily logic sheet at row 324,
# pushing the existing data (old rows 324-387) down to rows 326-389.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows starting at row 324 (shifts rows 324:387 -> 326:389)
$ws.Range("A324:R325").EntireRow.Insert()

# --- New row 324 ---
$ws.Cells.Item(324, 1).Value2 = 9
$ws.Cells.Item(324, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(324, 3).Value2 = "Metropolitana"
$ws.Cells.Item(324, 4).Value2 = 44711
$ws.Cells.Item(324, 5).Value2 = 13
$ws.Cells.Item(324, 6).Value2 = 100112032
$ws.Cells.Item(324, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(324, 8).Value2 = "Bola 8"
$ws.Cells.Item(324, 9).Value2 = "Primera"
$ws.Cells.Item(324, 10).Value2 = 34
$ws.Cells.Item(324, 11).Value2 = 18000
$ws.Cells.Item(324, 12).Value2 = 18000
$ws.Cells.Item(324, 13).Value2 = 18000
$ws.Cells.Item(324, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(324, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(324, 16).Value2 = 360
$ws.Cells.Item(324, 17).Value2 = 50
$ws.Cells.Item(324, 18).Value2 = "Hortaliza"

# --- New row 325 ---
$ws.Cells.Item(325, 1).Value2 = 9
$ws.Cells.Item(325, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(325, 3).Value2 = "Metropolitana"
$ws.Cells.Item(325, 4).Value2 = 44711
$ws.Cells.Item(325, 5).Value2 = 13
$ws.Cells.Item(325, 6).Value2 = 100112032
$ws.Cells.Item(325, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(325, 8).Value2 = "Sin especificar"
$ws.Cells.Item(325, 9).Value2 = "Primera"
$ws.Cells.Item(325, 10).Value2 = 70
$ws.Cells.Item(325, 11).Value2 = 15000
$ws.Cells.Item(325, 12).Value2 = 16000
$ws.Cells.Item(325, 13).Value2 = 15500
$ws.Cells.Item(325, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(325, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(325, 16).Value2 = 310
$ws.Cells.Item(325, 17).Value2 = 50
$ws.Cells.Item(325, 18).Value2 = "Hortaliza"
